$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H9").Formula = "=com.sun.star.sheet.addin.Analysis.getMround(G9,1)"
